$d = $word.ActiveDocument

# 1. Remove the stray "_GoBack" bookmark left over from the last edit session.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# 2. Remove the "For Spanish speakers..." paragraph entirely (it was an
#    end-of-document statement that is no longer needed).
$r = $d.Content
$found = $r.Find.Execute("For Spanish speakers*_________", $false, $false, $true, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $p = $r.Paragraphs(1)
    $p.Range.Delete()
}
